$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update/extend the data rows (row 1 = header, unchanged).
# Rows 2-3 get new values; rows 4-7 are newly added combinations
# of Sending cluster (ECs/FAPs/sCs) x Target cluster (FAPs/sCs).

# Row 2: ECs -> FAPs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efnb3"
$ws.Cells.Item(2, 3).Value = "Ephb6"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.1533166666666667
$ws.Cells.Item(2, 8).Value = 0.45995
$ws.Cells.Item(2, 9).Value = 0.1117088182569538
$ws.Cells.Item(2, 10).Value = 0.1117088182569538
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.7341896666666666
$ws.Cells.Item(2, 14).Value = 2.202569
$ws.Cells.Item(2, 15).Value = 0.4912907638668469
$ws.Cells.Item(2, 16).Value = 0.4912907638668469
$ws.Cells.Item(2, 17).Value = 0.1125635123944445
$ws.Cells.Item(2, 18).Value = 1.01307161155
$ws.Cells.Item(2, 19).Value = 0.05488151065212164
$ws.Cells.Item(2, 20).Value = 0.05488151065212162

# Row 3: ECs -> sCs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efnb3"
$ws.Cells.Item(3, 3).Value = "Ephb6"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.1533166666666667
$ws.Cells.Item(3, 8).Value = 0.45995
$ws.Cells.Item(3, 9).Value = 0.1117088182569538
$ws.Cells.Item(3, 10).Value = 0.1117088182569538
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.7602199999999999
$ws.Cells.Item(3, 14).Value = 2.28066
$ws.Cells.Item(3, 15).Value = 0.5087092361331531
$ws.Cells.Item(3, 16).Value = 0.5087092361331531
$ws.Cells.Item(3, 17).Value = 0.1165543963333333
$ws.Cells.Item(3, 18).Value = 1.048989567
$ws.Cells.Item(3, 19).Value = 0.05682730760483222
$ws.Cells.Item(3, 20).Value = 0.05682730760483221

# Row 4: FAPs -> FAPs
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Efnb3"
$ws.Cells.Item(4, 3).Value = "Ephb6"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.9559960000000002
$ws.Cells.Item(4, 8).Value = 2.867988
$ws.Cells.Item(4, 9).Value = 0.6965529954454279
$ws.Cells.Item(4, 10).Value = 0.6965529954454278
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.7341896666666666
$ws.Cells.Item(4, 14).Value = 2.202569
$ws.Cells.Item(4, 15).Value = 0.4912907638668469
$ws.Cells.Item(4, 16).Value = 0.4912907638668469
$ws.Cells.Item(4, 17).Value = 0.7018823845746668
$ws.Cells.Item(4, 18).Value = 6.316941461172001
$ws.Cells.Item(4, 19).Value = 0.3422100532061246
$ws.Cells.Item(4, 20).Value = 0.3422100532061246

# Row 5: FAPs -> sCs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efnb3"
$ws.Cells.Item(5, 3).Value = "Ephb6"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.9559960000000002
$ws.Cells.Item(5, 8).Value = 2.867988
$ws.Cells.Item(5, 9).Value = 0.6965529954454279
$ws.Cells.Item(5, 10).Value = 0.6965529954454278
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.7602199999999999
$ws.Cells.Item(5, 14).Value = 2.28066
$ws.Cells.Item(5, 15).Value = 0.5087092361331531
$ws.Cells.Item(5, 16).Value = 0.5087092361331531
$ws.Cells.Item(5, 17).Value = 0.7267672791200001
$ws.Cells.Item(5, 18).Value = 6.54090551208
$ws.Cells.Item(5, 19).Value = 0.3543429422393033
$ws.Cells.Item(5, 20).Value = 0.3543429422393032

# Row 6: sCs -> FAPs
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Efnb3"
$ws.Cells.Item(6, 3).Value = "Ephb6"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.2631543333333333
$ws.Cells.Item(6, 8).Value = 0.789463
$ws.Cells.Item(6, 9).Value = 0.1917381862976183
$ws.Cells.Item(6, 10).Value = 0.1917381862976183
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.7341896666666666
$ws.Cells.Item(6, 14).Value = 2.202569
$ws.Cells.Item(6, 15).Value = 0.4912907638668469
$ws.Cells.Item(6, 16).Value = 0.4912907638668469
$ws.Cells.Item(6, 17).Value = 0.1932051922718889
$ws.Cells.Item(6, 18).Value = 1.738846730447
$ws.Cells.Item(6, 19).Value = 0.0941992000086007
$ws.Cells.Item(6, 20).Value = 0.09419920000860069

# Row 7: sCs -> sCs
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Efnb3"
$ws.Cells.Item(7, 3).Value = "Ephb6"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.2631543333333333
$ws.Cells.Item(7, 8).Value = 0.789463
$ws.Cells.Item(7, 9).Value = 0.1917381862976183
$ws.Cells.Item(7, 10).Value = 0.1917381862976183
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.7602199999999999
$ws.Cells.Item(7, 14).Value = 2.28066
$ws.Cells.Item(7, 15).Value = 0.5087092361331531
$ws.Cells.Item(7, 16).Value = 0.5087092361331531
$ws.Cells.Item(7, 17).Value = 0.2000551872866666
$ws.Cells.Item(7, 18).Value = 1.80049668558
$ws.Cells.Item(7, 19).Value = 0.09753898628901761
$ws.Cells.Item(7, 20).Value = 0.09753898628901761
